$d = $word.ActiveDocument

# Locate the "Table 4.1" BodyText paragraph that immediately precedes the
# `import_table(...)` SourceCode chunk - the new R chunk
# (`EviewsR::eviews_commands("wfcreate m 1990 +90")`) must be inserted right
# after it and right before that SourceCode paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "Table 4.1`r") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Table 4.1' paragraph"
}

$anchorParagraph = $d.Paragraphs.Item($targetIndex)
$anchorParagraph.Range.InsertParagraphAfter() | Out-Null

# Re-fetch the freshly created (empty) paragraph and give it the SourceCode style.
$newParagraph = $d.Paragraphs.Item($targetIndex + 1)
$newParagraph.Style = "Source Code"

# Build the highlighted R chunk one syntax token at a time, matching the
# rStyle runs pandoc/highlighting would produce for:
#   EviewsR::eviews_commands("wfcreate m 1990 +90")
$segments = @(
    @("EviewsR", "NormalTok"),
    @("::", "SpecialCharTok"),
    @("eviews_commands", "FunctionTok"),
    @("(", "NormalTok"),
    @('"wfcreate m 1990 +90"', "StringTok"),
    @(")", "NormalTok")
)

$pos = $newParagraph.Range.Start
foreach ($seg in $segments) {
    $text = $seg[0]
    $styleName = $seg[1]

    $insertionPoint = $d.Range($pos, $pos)
    $insertionPoint.InsertAfter($text)

    $segEnd = $pos + $text.Length
    $segRange = $d.Range($pos, $segEnd)
    $segRange.Style = $styleName

    $pos = $segEnd
}
